$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1781.625
$ws.Range("J17").Value = 1781.625
$ws.Range("L17").Value = 5344.875
$ws.Range("N17").Value = -5680.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 48386.76
$ws.Range("I33").Value = 125494.625
$ws.Range("J33").Value = 935.7692
$ws.Range("K33").Value = 125494.625
$ws.Range("L33").Value = 935.7692
$ws.Range("M33").Value = -125265.625
$ws.Range("N33").Value = -1393.7692

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 922.2727
$ws.Range("I53").Value = 906.2727
$ws.Range("J53").Value = 954.2727
$ws.Range("K53").Value = 906.2727
$ws.Range("L53").Value = 954.2727
$ws.Range("M53").Value = -269.2727
$ws.Range("N53").Value = -2228.2727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4773.778
$ws.Range("I64").Value = 4373.727
$ws.Range("K64").Value = 4373.727
$ws.Range("M64").Value = -4125.727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4773.778
$ws.Range("I67").Value = 4373.727
$ws.Range("K67").Value = 4373.727
$ws.Range("M67").Value = -3515.727

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4199.4
$ws.Range("I76").Value = 4666
$ws.Range("K76").Value = 4666
$ws.Range("M76").Value = -4351

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4199.4
$ws.Range("I79").Value = 4666
$ws.Range("K79").Value = 4666
$ws.Range("M79").Value = -3574

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 281.22726
$ws.Range("I92").Value = 264
$ws.Range("J92").Value = 311.375
$ws.Range("K92").Value = 264
$ws.Range("L92").Value = 311.375
$ws.Range("M92").Value = 984
$ws.Range("N92").Value = -2807.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1340.4546
$ws.Range("I106").Value = 1224.5
$ws.Range("K106").Value = 1224.5
$ws.Range("M106").Value = -593.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2576.6287
$ws.Range("I138").Value = 1727.0588
$ws.Range("J138").Value = 3379
$ws.Range("K138").Value = 5181.1764
$ws.Range("L138").Value = 10137
$ws.Range("M138").Value = -41.17640000000029
$ws.Range("N138").Value = -20417

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 814.386
$ws.Range("I32").Value = 814.386
$ws.Range("K32").Value = 814.386
$ws.Range("M32").Value = -527.386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2535.879
$ws.Range("I61").Value = 1820.0454
$ws.Range("K61").Value = 1820.0454
$ws.Range("M61").Value = -1608.0454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4210758.5
$ws.Range("I74").Value = 2647367.2
$ws.Range("J74").Value = 6946693
$ws.Range("K74").Value = 2647367.2
$ws.Range("L74").Value = 6946693
$ws.Range("M74").Value = -2646493.2
$ws.Range("N74").Value = -6948441

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4210758.5
$ws.Range("I77").Value = 2647367.2
$ws.Range("J77").Value = 6946693
$ws.Range("K77").Value = 13236836
$ws.Range("L77").Value = 34733465
$ws.Range("M77").Value = -13232468
$ws.Range("N77").Value = -34742201

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 737
$ws.Range("I102").Value = 699.7143
$ws.Range("K102").Value = 699.7143
$ws.Range("M102").Value = 922.2857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2535.879
$ws.Range("I136").Value = 1820.0454
$ws.Range("K136").Value = 5460.1362
$ws.Range("M136").Value = -2910.1362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3296.8333
$ws.Range("I86").Value = 2318.5
$ws.Range("J86").Value = 5253.5
$ws.Range("K86").Value = 2318.5
$ws.Range("L86").Value = 5253.5
$ws.Range("M86").Value = -1195.5
$ws.Range("N86").Value = -7499.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3296.8333
$ws.Range("I89").Value = 2318.5
$ws.Range("J89").Value = 5253.5
$ws.Range("K89").Value = 11592.5
$ws.Range("L89").Value = 26267.5
$ws.Range("M89").Value = -5976.5
$ws.Range("N89").Value = -37499.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2895.6667
$ws.Range("I99").Value = 2927.5
$ws.Range("K99").Value = 2927.5
$ws.Range("M99").Value = -1429.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3588.4
$ws.Range("J105").Value = 4800
$ws.Range("L105").Value = 4800
$ws.Range("N105").Value = -8294

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2464.3157
$ws.Range("I107").Value = 909.2308
$ws.Range("K107").Value = 909.2308
$ws.Range("M107").Value = 1010.7692

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 83347250
$ws.Range("I134").Value = 41682896
$ws.Range("K134").Value = 125048688
$ws.Range("M134").Value = -125046153

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2213.5715
$ws.Range("J16").Value = 2643.3333
$ws.Range("L16").Value = 2643.3333
$ws.Range("N16").Value = -3217.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1906.3572
$ws.Range("I58").Value = 1976.0769
$ws.Range("K58").Value = 1976.0769
$ws.Range("M58").Value = -1773.0769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2650.0908
$ws.Range("I105").Value = 2158.0667
$ws.Range("K105").Value = 2158.0667
$ws.Range("M105").Value = -411.0666999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 5976.9565
$ws.Range("I107").Value = 6999.647
$ws.Range("K107").Value = 6999.647
$ws.Range("M107").Value = -5079.647

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2213.5715
$ws.Range("J113").Value = 2643.3333
$ws.Range("L113").Value = 2643.3333
$ws.Range("N113").Value = -6983.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1719.6111
$ws.Range("J122").Value = 2119.7273
$ws.Range("L122").Value = 6359.1819
$ws.Range("N122").Value = -11259.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1906.3572
$ws.Range("I136").Value = 1976.0769
$ws.Range("K136").Value = 5928.2307
$ws.Range("M136").Value = -3378.2307

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1839500
$ws.Range("I46").Value = 5500000
$ws.Range("J46").Value = 9250
$ws.Range("K46").Value = 16500000
$ws.Range("L46").Value = 27750
$ws.Range("M46").Value = -16499909
$ws.Range("N46").Value = -27932

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 665.6667
$ws.Range("J50").Value = 1500
$ws.Range("L50").Value = 4500
$ws.Range("N50").Value = -5462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 665.6667
$ws.Range("J53").Value = 1500
$ws.Range("L53").Value = 4500
$ws.Range("N53").Value = -5462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1840.8572
$ws.Range("J129").Value = 2573.5
$ws.Range("L129").Value = 7720.5
$ws.Range("N129").Value = -17720.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4299.3335
$ws.Range("I80").Value = 3085
$ws.Range("K80").Value = 3085
$ws.Range("M80").Value = -2087

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4299.3335
$ws.Range("I83").Value = 3085
$ws.Range("K83").Value = 15425
$ws.Range("M83").Value = -10433

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1170.0526
$ws.Range("I113").Value = 1184.9286
$ws.Range("J113").Value = 1128.4
$ws.Range("K113").Value = 1184.9286
$ws.Range("L113").Value = 1128.4
$ws.Range("M113").Value = 985.0714
$ws.Range("N113").Value = -5468.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3474.4736
$ws.Range("I122").Value = 3712.5715
$ws.Range("K122").Value = 11137.7145
$ws.Range("M122").Value = -8687.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1995
$ws.Range("I132").Value = 1995
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5985
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3455
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4999
$ws.Range("I61").Value = 4999
$ws.Range("K61").Value = 4999
$ws.Range("M61").Value = -4797

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2069.5557
$ws.Range("I93").Value = 2132.2856
$ws.Range("K93").Value = 2132.2856
$ws.Range("M93").Value = -884.2856000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4999
$ws.Range("I113").Value = 4999
$ws.Range("K113").Value = 4999
$ws.Range("M113").Value = -2829

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2192.3333
$ws.Range("I100").Value = 1787.75
$ws.Range("K100").Value = 3575.5
$ws.Range("M100").Value = -3034.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1760.5
$ws.Range("I107").Value = 791.1111
$ws.Range("K107").Value = 2373.3333
$ws.Range("M107").Value = -453.3332999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1929.2222
$ws.Range("I132").Value = 1115.2
$ws.Range("K132").Value = 3345.6
$ws.Range("M132").Value = -815.6000000000004
